$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.798.34"
$ws.Range("E2").Value = "  +0.58%  "

$ws.Range("D3").Value = "1.703.26"
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  +0.85%  "

$ws.Range("D5").Value = "316.40"
$ws.Range("E5").Value = "  -0.09%  "

$ws.Range("E6").Value = "  +0.75%  "

$ws.Range("D7").Value = "0.3943"
$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").Value = "0.4075"
$ws.Range("E8").Value = "  +0.60%  "

$ws.Range("D9").Value = "1.516"
$ws.Range("E9").Value = "  +1.87%  "

$ws.Range("D10").Value = "1.012"
$ws.Range("E10").Value = "  +0.98%  "

$ws.Range("D11").Value = "52.47"
$ws.Range("E11").Value = "  +1.14%  "

$ws.Range("D12").Value = "0.08810"
$ws.Range("E12").Value = "  -0.76%  "

$ws.Range("D13").Value = "7.615"
$ws.Range("E13").Value = "  +6.30%  "

$ws.Range("D14").Value = "24.76"
$ws.Range("E14").Value = "  +5.45%  "

$ws.Range("D15").Value = "0.00001372"
$ws.Range("E15").Value = "  +3.53%  "

$ws.Range("D16").Value = "8.046"
$ws.Range("E16").Value = "  -1.17%  "

$ws.Range("D17").Value = "1.704.91"
$ws.Range("E17").Value = "  +0.66%  "

$ws.Range("D18").Value = "99.52"
$ws.Range("E18").Value = "  -0.48%  "

$ws.Range("D19").Value = "0.07120"
$ws.Range("E19").Value = "  +1.54%  "

$ws.Range("D20").Value = "19.98"
$ws.Range("E20").Value = "  +1.70%  "

$ws.Range("D21").Value = "7.400"
$ws.Range("E21").Value = "  +5.49%  "

$ws.Range("D22").Value = "1.011"
$ws.Range("E22").Value = "  +1.03%  "

$ws.Range("D23").Value = "14.38"
$ws.Range("E23").Value = "  +0.27%  "

$ws.Range("D24").Value = "24.791.73"
$ws.Range("E24").Value = "  +0.54%  "

$ws.Range("D25").Value = "3.052"
$ws.Range("E25").Value = "  -4.71%  "

$ws.Range("D26").Value = "2.349"
$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("D27").Value = "22.82"
$ws.Range("E27").Value = "  +0.55%  "

$ws.Range("D28").Value = "165.01"
$ws.Range("E28").Value = "  +0.98%  "

$ws.Range("D29").Value = "8.717"
$ws.Range("E29").Value = "  +16.42%  "

$ws.Range("D30").Value = "138.45"
$ws.Range("E30").Value = "  +1.42%  "

$ws.Range("D31").Value = "5.225"
$ws.Range("E31").Value = "  +1.16%  "

$ws.Range("D32").Value = "7.721"
$ws.Range("E32").Value = "  +7.89%  "

$ws.Range("D33").Value = "1.891.96"
$ws.Range("E33").Value = "  +0.58%  "

$ws.Range("D34").Value = "0.08913"
$ws.Range("E34").Value = "  +3.82%  "

$ws.Range("D35").Value = "1.050"
$ws.Range("E35").Value = "  -1.63%  "

$ws.Range("D36").Value = "1.987"
$ws.Range("E36").Value = "  +3.25%  "

$ws.Range("D37").Value = "0.2751"
$ws.Range("E37").Value = "  +0.49%  "

$ws.Range("D38").Value = "0.02896"
$ws.Range("E38").Value = "  +6.27%  "

$ws.Range("D39").Value = "10.84"
$ws.Range("E39").Value = "  -5.32%  "

$ws.Range("D40").Value = "14.43"
$ws.Range("E40").Value = "  +0.02%  "

$ws.Range("D41").Value = "0.09155"
$ws.Range("E41").Value = "  -0.12%  "

$ws.Range("D42").Value = "0.7876"
$ws.Range("E42").Value = "  +2.75%  "

$ws.Range("D43").Value = "1.476"
$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("D44").Value = "16.74"
$ws.Range("E44").Value = "  +4.26%  "

$ws.Range("D45").Value = "0.7255"
$ws.Range("E45").Value = "  +1.03%  "

$ws.Range("D46").Value = "2.613"
$ws.Range("E46").Value = "  +0.66%  "

$ws.Range("D47").Value = "4.251"
$ws.Range("E47").Value = "  +0.76%  "

$ws.Range("E48").Value = "  +0.73%  "

$ws.Range("D49").Value = "1.333"
$ws.Range("E49").Value = "  +0.24%  "

$ws.Range("D50").Value = "140.41"
$ws.Range("E50").Value = "  +0.00%  "

$ws.Range("D51").Value = "91.98"
$ws.Range("E51").Value = "  +2.57%  "
